$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(59, 'diary', '日記|にっき')
    ,@(60, 'to fill in', '記入する|きにゅうする')
    ,@(61, '(news) article', '記事|きじ')
    ,@(62, 'to memorize', '暗記する|あんきする')
    ,@(63, 'bank', '銀行|ぎんこう')
    ,@(64, 'silver medal', '銀メダル|ぎんメダル')
    ,@(65, 'land covered with snow', '銀世界|ぎんせかい')
    ,@(66, 'one time', '一回|いっかい')
    ,@(67, 'out-of-service bus', '回送バス|かいそうバス')
    ,@(68, 'last inning; last episode', '最終回|さいしゅうかい')
    ,@(69, 'to turn', '回す|まわす')
    ,@(70, 'evening', '夕方|ゆうがた')
    ,@(71, 'dinner', '夕食|ゆうしょく')
    ,@(72, 'Tanabata', '七夕|たなばた')
    ,@(73, 'setting sun', '夕日|ゆうひ')
    ,@(74, 'evening newspaper', '夕刊|ゆうかん')
    ,@(75, 'Mr./Ms. Kuroki', '黒木さん|くろきさん')
    ,@(76, 'black', '黒い|くろい')
    ,@(77, 'black and white photograph', '白黒写真|しろくろしゃしん')
    ,@(78, 'blackboard', '黒板|こくばん')
    ,@(79, 'a thing to take care of', '用事|ようじ')
    ,@(80, 'to prepare', '用意する|よういする')
    ,@(81, 'for children', '子供用|こどもよう')
    ,@(82, 'cost', '費用|ひよう')
    ,@(83, 'weekend', '週末|しゅうまつ')
    ,@(84, 'end of the month', '月末|げつまつ')
    ,@(85, 'year-end', '年末|ねんまつ')
    ,@(86, 'final examination', '期末試験|きまつしけん')
    ,@(87, 'the end', '末|すえ')
    ,@(88, 'to wait', '待つ|まつ')
    ,@(89, 'waiting room', '待合室|まちあいしつ')
    ,@(90, 'to expect', '期待する|きたいする')
    ,@(91, 'invitation', '招待|しょうたい')
    ,@(92, 'over-time work', '残業|ざんぎょう')
    ,@(93, 'to leave', '残す|のこす')
    ,@(94, 'regrettable', '残念|ざんねん')
    ,@(95, 'regret', '心残り|こころのこり')
    ,@(96, 'account balance', '残高|ざんだか')
    ,@(97, 'station', '駅|えき')
    ,@(98, 'Tokyo Station', '東京駅|とうきょうえき')
    ,@(99, 'train station attendant', '駅員|えきいん')
    ,@(100, 'near/in front of the station', '駅前|えきまえ')
    ,@(101, 'the first', '一番|いちばん')
    ,@(102, 'number', '番号|ばんごう')
    ,@(103, 'TV program', '番組|ばんぐみ')
    ,@(104, 'police box', '交番|こうばん')
    ,@(105, 'turn; order', '順番|じゅんばん')
    ,@(106, 'to explain', '説明する|せつめいする')
    ,@(107, 'novel', '小説|しょうせつ')
    ,@(108, 'novelist', '小説家|しょうせつか')
    ,@(109, 'to preach', '説教する|せっきょうする')
    ,@(110, 'to guide', '案内する|あんないする')
    ,@(111, 'information desk', '案内所|あんないじょ')
    ,@(112, 'idea; proposal', '案|あん')
    ,@(113, 'proposal', '提案|ていあん')
    ,@(114, 'my wife', '家内|かない')
    ,@(115, 'domestic', '国内|こくない')
    ,@(116, 'internal medicine', '内科|ないか')
    ,@(117, 'inside', '内側|うちがわ')
    ,@(118, 'to forget', '忘れる|わすれる')
    ,@(119, 'lost article', '忘れ物|わすれもの')
    ,@(120, 'year-end party', '忘年会|ぼうねんかい')
    ,@(121, 'to keep (a promise)', '守る|まもる')
    ,@(122, 'absence; not at home', '留守|るす')
    ,@(123, 'answering machine', '留守番電話|るすばんでんわ')
    ,@(124, 'charm', 'お守り|おまもり')
    ,@(125, 'security guard', '守衛|しゅえい')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
